$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Cover sheet updates
# ---------------------------------------------------------------------------
$cover = $wb.Worksheets.Item("Cover")

# A3: components list - update VistALink build timestamp (084015 -> 085649)
# This is a multi-run rich string; rebuild the full text then re-apply the
# per-run font formatting (Arial / Calibri alternating, bold, size 12).
$r = $cover.Range("A3")
$em = [char]0x2014
$full = "1.  Web UI $em TAS.02.00.4_20190502_140102" + "`n" `
      + "2.  TAS API $em TAS_API.01.00.65_20190422_212745" + "`n" `
      + "3.  FHIR $em VA-FHIR-Server_02.00.5_20190507_141937" + "`n" `
      + "4.  VistALink $em VAVLGateway_02.00.3_20190516_085649"
$r.Value = $full

$runs = @(
  @(0,11,"Arial"),
  @(11,1,"Calibri"),
  @(12,41,"Arial"),
  @(53,1,"Calibri"),
  @(54,43,"Arial"),
  @(97,1,"Calibri"),
  @(98,54,"Arial"),
  @(152,1,"Calibri"),
  @(153,36,"Arial")
)
foreach ($run in $runs) {
  $c = $r.Characters($run[0]+1, $run[1])
  $c.Font.Name = $run[2]
  $c.Font.Bold = $true
  $c.Font.Size = 12
}

# A15: revision date shown on cover -> May 20, 2019
$cover.Range("A15").Value = "May 20, 2019"

# A16: version shown on cover -> Version 5.0
$cover.Range("A16").Value = "Version 5.0"

# ---------------------------------------------------------------------------
# 2) Revision History sheet updates
# ---------------------------------------------------------------------------
$rh = $wb.Worksheets.Item("Revision History")

# Insert a new row above the current row 3 for the Build 10 / v5.0 entry
$rh.Rows.Item(3).Insert()

# Copy the formatting of the (now shifted) row 4 back up onto new row 3
$rh.Range("A4:D4").Copy()
$rh.Range("A3:D3").PasteSpecial(-4122)
$rh.Rows.Item(3).RowHeight = $rh.Rows.Item(4).RowHeight()

# Fill in the new revision-history row values
$rh.Range("A3").Value = "05-20-2019"
$rh.Range("B3").Value = "5.0"
$rh.Range("C3").Value = "Updated VistALink version number "
$rh.Range("D3").Value = "Donald Fowlds"

# ---------------------------------------------------------------------------
# 3) Selection / active cell bookkeeping to match the authored workbook
# ---------------------------------------------------------------------------
$cover.Activate()
$cover.Range("A2").Select()

$rh.Activate()
$rh.Range("A3").Select()

$cover.Activate()
